# Adds two new "base"/"external" system commands to the hidden '#system'
# sheet of the nexial-core unit-test workbook:
#   - base:     assertMatch(text,regex)   -> inserted alphabetically into column F
#   - external: openFile(filePath)        -> inserted alphabetically into column J
#
# Inserting a row into the middle of a column-based lookup list (without
# touching the other lists that live in neighbouring columns) also causes a
# cascade of changes that are visible in the target OOXML diff:
#   - the "target" named range (column A) used to duplicate the "tn.5250"
#     command list; that duplicate entry is dropped, shifting the remaining
#     entries (web/webalert/webcookie/ws/ws.async/xml) up by one row.
#   - columns AA..AG shift left by one column (AA's original "tn.5250" data
#     is superseded by what used to be in AB, etc.), leaving column AG empty.
#
# All of this is reproduced below with plain cell-by-cell Value copies so
# that only the intended columns/rows move (Range.Insert affects whole
# rows/columns in this host, which would be too broad).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) column F ("base"): insert "assertMatch(text,regex)" at F11, pushing
#    the existing F11:F44 block down to F12:F45.
# ---------------------------------------------------------------------
for ($r = 44; $r -ge 11; $r--) {
    $val = $ws.Cells.Item($r, 6).Value()
    $ws.Cells.Item($r + 1, 6).Value = $val
}
$ws.Cells.Item(11, 6).Value = "assertMatch(text,regex)"

# ---------------------------------------------------------------------
# 2) column J ("external"): insert "openFile(filePath)" at J2, pushing
#    the existing J2:J6 block down to J3:J7.
# ---------------------------------------------------------------------
for ($r = 6; $r -ge 2; $r--) {
    $val = $ws.Cells.Item($r, 10).Value()
    $ws.Cells.Item($r + 1, 10).Value = $val
}
$ws.Cells.Item(2, 10).Value = "openFile(filePath)"

# ---------------------------------------------------------------------
# 3) column A ("target"): drop the "tn.5250" duplicate at A27, pulling
#    A28:A33 up to A27:A32 and clearing the now-unused A33.
# ---------------------------------------------------------------------
for ($r = 27; $r -le 32; $r++) {
    $val = $ws.Cells.Item($r + 1, 1).Value()
    $ws.Cells.Item($r, 1).Value = $val
}
$ws.Cells.Item(33, 1).ClearContents()

# ---------------------------------------------------------------------
# 4) columns AA..AG: shift each column one step to the left (AB->AA,
#    AC->AB, AD->AC, AE->AD, AF->AE, AG->AF), clearing column AG
#    afterwards since nothing shifts into it.
# ---------------------------------------------------------------------
for ($col = 27; $col -le 32; $col++) {
    for ($r = 1; $r -le 151; $r++) {
        $val = $ws.Cells.Item($r, $col + 1).Value()
        $ws.Cells.Item($r, $col).Value = $val
    }
}
$ws.Range("AG1:AG151").ClearContents()

# ---------------------------------------------------------------------
# 5) update the defined names (Names collection) to match the new
#    extents produced by the row/column shifts above.
# ---------------------------------------------------------------------
$newRefersTo = @{
    "base"     = "='#system'!`$F`$2:`$F`$45"
    "external" = "='#system'!`$J`$2:`$J`$7"
    "target"   = "='#system'!`$A`$2:`$A`$32"
    "web"      = "='#system'!`$AA`$2:`$AA`$151"
    "webalert" = "='#system'!`$AB`$2:`$AB`$8"
    "webcookie"= "='#system'!`$AC`$2:`$AC`$10"
    "ws"       = "='#system'!`$AD`$2:`$AD`$17"
    "ws.async" = "='#system'!`$AE`$2:`$AE`$8"
    "xml"      = "='#system'!`$AF`$2:`$AF`$27"
}

$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $item = $names.Item($i)
    $nm = $item.Name()
    if ($newRefersTo.ContainsKey($nm)) {
        $item.RefersTo = $newRefersTo[$nm]
    }
}

Write-Host "done"
